$d = $word.ActiveDocument

$replacements = @(
    @{old = "199÷7=28, 3"; new = "756÷6=126, 0"},
    @{old = "514÷5=102, 4"; new = "691÷6=115, 1"},
    @{old = "486÷9=54, 0"; new = "894÷5=178, 4"},
    @{old = "878÷4=219, 2"; new = "620÷4=155, 0"},
    @{old = "275÷2=137, 1"; new = "629÷6=104, 5"},
    @{old = "946÷8=118, 2"; new = "432÷9=48, 0"},
    @{old = "623÷2=311, 1"; new = "610÷5=122, 0"},
    @{old = "377÷7=53, 6"; new = "226÷5=45, 1"},
    @{old = "810÷9=90, 0"; new = "197÷4=49, 1"},
    @{old = "794÷3=264, 2"; new = "638÷6=106, 2"},
    @{old = "337÷8=42, 1"; new = "597÷6=99, 3"},
    @{old = "461÷6=76, 5"; new = "546÷8=68, 2"},
    @{old = "201÷2=100, 1"; new = "325÷3=108, 1"},
    @{old = "195÷3=65, 0"; new = "867÷6=144, 3"},
    @{old = "558÷7=79, 5"; new = "705÷4=176, 1"},
    @{old = "232÷9=25, 7"; new = "302÷5=60, 2"},
    @{old = "898÷4=224, 2"; new = "347÷6=57, 5"},
    @{old = "405÷7=57, 6"; new = "448÷4=112, 0"},
    @{old = "545÷2=272, 1"; new = "931÷8=116, 3"},
    @{old = "786÷8=98, 2"; new = "834÷4=208, 2"},
    @{old = "887÷3=295, 2"; new = "729÷2=364, 1"},
    @{old = "450÷4=112, 2"; new = "351÷3=117, 0"},
    @{old = "822÷7=117, 3"; new = "549÷3=183, 0"},
    @{old = "756÷9=84, 0"; new = "857÷3=285, 2"},
    @{old = "147÷7=21, 0"; new = "769÷2=384, 1"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
